# Edit slide 9 ("Neural Network") content placeholder text.
# Apply targeted substring replacements (via TextRange.Characters(start,length))
# instead of rewriting the whole TextRange.Text, so that paragraph breaks,
# run boundaries/properties (e.g. the "Keras" runs' err="1" spellcheck flag)
# and the trailing empty paragraphs are preserved exactly as in the original.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Apply edits from the end of the text backwards so that earlier character
# offsets stay valid even though some replacements change the text length.

$tr.Characters(438, 27).Text = " result was slightly better."
$tr.Characters(369, 64).Text = " prediction results with the SVM result under ROC. It showed the "
$tr.Characters(352, 12).Text = "Compared the "
$tr.Characters(165, 186).Text = "Created a Sequential NN with 5 layers (input/output + 3 hidden layers). Another feature we added was 'dropout' which avoided overfitting and reduced running time."
$tr.Characters(1, 163).Text = "Created a plotting function to plot the loss/accuracy against the number of epochs on the training and validation dataset. This works similar as tensor board."
